# Refactor, hamcrest and division of features
# Adds a new "Movies" worksheet (with a single "name" column containing
# "Fight Club") after the existing "Filter" sheet, and makes it the active
# / selected sheet.

$wb = $excel.ActiveWorkbook

# Excel inserts new sheets before the active sheet, so add it and then
# move it to the end of the tab strip.
$ws = $wb.Worksheets.Add()
$ws.Name = "Movies"
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-resolve the sheet by name: Move() invalidates the old COM reference.
$movies = $wb.Worksheets.Item("Movies")

$movies.Range("A1").Value = "name"
$movies.Range("A2").Value = "Fight Club"

$movies.Columns.Item(1).ColumnWidth = 10.33

# Make "Movies" the active/selected sheet and park the selection on P21,
# matching the author's saved view state.
$movies.Select() | Out-Null
$movies.Range("P21").Select() | Out-Null
